# Update the "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to match the latest scraped counts.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 7084
    4  = 4543
    5  = 68
    9  = 95
    10 = 67
    11 = 63
    13 = 617
    14 = 123
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
